# Commit message: "Added the optional skip_nodes parameter to skip desired nodes"
#
# The underlying workbook edit adds a third worksheet ("Sheet1") holding a
# reduced subset of the TSP points found on "raw" (nodes 4-17 dropped,
# keeping nodes 1,2,3,18-25), wires up a matching sheet-scoped "tsp"
# defined name for it, and leaves a couple of small selection / page-setup
# breadcrumbs behind on the other two sheets from the edit session.

$wb = $excel.ActiveWorkbook

$wsSimplified = $wb.Worksheets.Item("simplified")
$wsRaw = $wb.Worksheets.Item("raw")

# -- raw: switch the page to portrait and reselect the full data range -----
$wsRaw.PageSetup.Orientation = 1
$wsRaw.Range("A1:B25").Select() | Out-Null

# -- simplified: selection moves off the old A1:B21 block, down to D6 ------
$wsSimplified.Range("D6").Select() | Out-Null

# -- add the new "Sheet1" with the reduced / skipped point list ------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sheet1"

$points = @(
    @(1,  20833.333299999998, 17100),
    @(2,  20900,              17066.666700000002),
    @(3,  21300,              13016.6667),
    @(18, 26433.333299999998, 13433.3333),
    @(19, 26550,              13850),
    @(20, 26733.333299999998, 11683.3333),
    @(21, 27026.111099999998, 13051.9444),
    @(22, 27096.111099999998, 13415.8333),
    @(23, 27153.611099999998, 13203.3333),
    @(24, 27166.666700000002, 9833.3333000000002),
    @(25, 27233.333299999998, 10450)
)

for ($i = 0; $i -lt $points.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $points[$i][0]
    $ws.Cells.Item($r, 2).Value = $points[$i][1]
    $ws.Cells.Item($r, 3).Value = $points[$i][2]
}

# Sheet-scoped "tsp" name for the new sheet, matching the pattern already
# used by "simplified" and "raw".
$ws.Names.Add("tsp", "=Sheet1!`$B`$1:`$C`$11") | Out-Null

# Leave the selection where the edit session left it, just past the data.
$ws.Range("A12").Select() | Out-Null

Write-Output "done"
